$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.05"
$ws.Range("E2").Value = "'0.92%"
$ws.Range("D3").Value = "'29.62"
$ws.Range("E3").Value = "'-1.37%"
$ws.Range("D4").Value = "'5.153"
$ws.Range("E4").Value = "'-0.06%"
$ws.Range("D5").Value = "'0.05800"
$ws.Range("E5").Value = "'2.25%"
$ws.Range("D6").Value = "'6.653"
$ws.Range("E6").Value = "'1.70%"
$ws.Range("D7").Value = "'3.199"
$ws.Range("E7").Value = "'5.46%"
$ws.Range("D8").Value = "'0.8511"
$ws.Range("E8").Value = "'0.36%"
$ws.Range("D9").Value = "'0.8618"
$ws.Range("E9").Value = "'0.20%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("E10").Value = "'2.52%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07088"
$ws.Range("E11").Value = "'2.62%"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03366"
$ws.Range("E12").Value = "'0.63%"
$ws.Range("D13").Value = "'0.03198"
$ws.Range("E13").Value = "'10.67%"
$ws.Range("D14").Value = "'0.09375"
$ws.Range("E14").Value = "'0.03%"
$ws.Range("D15").Value = "'0.001522"
$ws.Range("E15").Value = "'-0.32%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005974"
$ws.Range("E16").Value = "'-0.07%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006152"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.487"
$ws.Range("E18").Value = "'-0.59%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.223"
$ws.Range("E19").Value = "'4.22%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3197"
$ws.Range("E20").Value = "'1.50%"
$ws.Range("E21").Value = "'-1.62%"
$ws.Range("D22").Value = "'3.477"
$ws.Range("E22").Value = "'-4.10%"
$ws.Range("D23").Value = "'0.04148"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("E24").Value = "'0.48%"
$ws.Range("E25").Value = "'1.07%"
$ws.Range("E26").Value = "'-6.67%"
$ws.Range("E27").Value = "'2.51%"
$ws.Range("E28").Value = "'3.79%"
$ws.Range("D40").Value = "'0.03745"
$ws.Range("E40").Value = "'-0.76%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1072"
$ws.Range("E41").Value = "'1.33%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002449"
$ws.Range("E42").Value = "'6.99%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003518"
$ws.Range("E43").Value = "'-39.94%"
$ws.Range("D44").Value = "'0.009184"
$ws.Range("E44").Value = "'-0.94%"
$ws.Range("D45").Value = "'0.00005288"
$ws.Range("E45").Value = "'3.64%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("D47").Value = "'0.05793"
$ws.Range("E47").Value = "'-35.61%"
$ws.Range("D48").Value = "'0.002174"
$ws.Range("E48").Value = "'-21.49%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.08%"
